$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.632.56"
$ws.Range("E2").Value = "  +2.62%  "

$ws.Range("D3").Value = "2.946.12"
$ws.Range("E3").Value = "  +2.23%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'590.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("D6").Value = "'147.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.53%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").Value = "2.945.24"
$ws.Range("E8").Value = "  +2.24%  "

$ws.Range("E9").Value = "  +3.31%  "

$ws.Range("D10").Value = "'7.19"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.67%  "

$ws.Range("E11").Value = "  +10.21%  "

$ws.Range("D12").Value = "'0.436"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.05%  "

$ws.Range("D13").Value = "'0.0000235"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.33%  "

$ws.Range("D14").Value = "'32.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.62%  "

$ws.Range("E15").Value = "  -0.68%  "

$ws.Range("D16").Value = "3.432.27"
$ws.Range("E16").Value = "  +2.10%  "

$ws.Range("D17").Value = "62.579.03"
$ws.Range("E17").Value = "  +2.61%  "

$ws.Range("E18").Value = "  +2.70%  "

$ws.Range("D19").Value = "2.948.29"
$ws.Range("E19").Value = "  +2.24%  "

$ws.Range("D20").Value = "'434.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.18%  "

$ws.Range("D21").Value = "'13.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.74%  "

$ws.Range("D22").Value = "'0.665"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.98%  "

$ws.Range("D23").Value = "'6.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.26%  "

$ws.Range("D24").Value = "'11.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.99%  "

$ws.Range("D25").Value = "'80.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.61%  "

$ws.Range("D26").Value = "'11.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.86%  "

$ws.Range("D27").Value = "'2.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.40%  "

$ws.Range("E28").Value = "  +0.07%  "

$ws.Range("D29").Value = "'7.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.39%  "

$ws.Range("D30").Value = "'2.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.78%  "

$ws.Range("D31").Value = "'2.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.25%  "

$ws.Range("D32").Value = "'0.0000102"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +21.00%  "

$ws.Range("E33").Value = "  +5.35%  "

$ws.Range("D34").Value = "'26.21"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.13%  "

$ws.Range("E35").Value = "  -0.05%  "

$ws.Range("D36").Value = "'0.990"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.42%  "

$ws.Range("E37").Value = "  +3.06%  "

$ws.Range("D38").Value = "'3.03"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.65%  "

$ws.Range("D39").Value = "'49.70"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.48%  "

$ws.Range("E40").Value = "  +6.28%  "

$ws.Range("D41").Value = "'8.37"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.80%  "

$ws.Range("E42").Value = "  +0.08%  "

$ws.Range("D43").Value = "'0.276"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.87%  "

$ws.Range("D44").Value = "'39.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.12%  "

$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.698.11"
$ws.Range("E45").Value = "  +1.52%  "

$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "'135.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.43%  "

$ws.Range("D47").Value = "'0.0340"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.55%  "

$ws.Range("D48").Value = "'353.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.26%  "

$ws.Range("E49").Value = "  +0.02%  "

$ws.Range("E50").Value = "  +2.27%  "

$ws.Range("D51").Value = "'22.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.81%  "
